# Refresh the Icam2-Itgam LR-pairs sheet with the recomputed TPM values.
# The new data drops the "-> ECs" target rows (NATMI no longer reports
# signalling onto ECs in this run) and the per-pair statistics for the
# remaining 6 sender/target combinations were recomputed with the new TPM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam2"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.31228666666667
$ws.Range("H2").Value = 63.93686
$ws.Range("I2").Value = 0.9506775731819035
$ws.Range("J2").Value = 0.9506775731819034
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005673666666666667
$ws.Range("N2").Value = 0.017021
$ws.Range("O2").Value = 0.1234828534325781
$ws.Range("P2").Value = 0.1234828534325781
$ws.Range("Q2").Value = 0.1209188104511111
$ws.Range("R2").Value = 1.08826929406
$ws.Range("S2").Value = 0.11739237943086
$ws.Range("T2").Value = 0.11739237943086

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam2"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.31228666666667
$ws.Range("H3").Value = 63.93686
$ws.Range("I3").Value = 0.9506775731819035
$ws.Range("J3").Value = 0.9506775731819034
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04027333333333333
$ws.Range("N3").Value = 0.12082
$ws.Range("O3").Value = 0.8765171465674219
$ws.Range("P3").Value = 0.876517146567422
$ws.Range("Q3").Value = 0.8583168250222223
$ws.Range("R3").Value = 7.7248514252
$ws.Range("S3").Value = 0.8332851937510434
$ws.Range("T3").Value = 0.8332851937510434

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icam2"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7500946666666666
$ws.Range("H4").Value = 2.250284
$ws.Range("I4").Value = 0.03345948693899053
$ws.Range("J4").Value = 0.03345948693899053
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005673666666666667
$ws.Range("N4").Value = 0.017021
$ws.Range("O4").Value = 0.1234828534325781
$ws.Range("P4").Value = 0.1234828534325781
$ws.Range("Q4").Value = 0.004255787107111111
$ws.Range("R4").Value = 0.038302083964
$ws.Range("S4").Value = 0.004131672921616629
$ws.Range("T4").Value = 0.00413167292161663

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam2"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7500946666666666
$ws.Range("H5").Value = 2.250284
$ws.Range("I5").Value = 0.03345948693899053
$ws.Range("J5").Value = 0.03345948693899053
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04027333333333333
$ws.Range("N5").Value = 0.12082
$ws.Range("O5").Value = 0.8765171465674219
$ws.Range("P5").Value = 0.876517146567422
$ws.Range("Q5").Value = 0.03020881254222222
$ws.Range("R5").Value = 0.2718793128799999
$ws.Range("S5").Value = 0.0293278140173739
$ws.Range("T5").Value = 0.0293278140173739

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Icam2"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3556153333333333
$ws.Range("H6").Value = 1.066846
$ws.Range("I6").Value = 0.01586293987910606
$ws.Range("J6").Value = 0.01586293987910605
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005673666666666667
$ws.Range("N6").Value = 0.017021
$ws.Range("O6").Value = 0.1234828534325781
$ws.Range("P6").Value = 0.1234828534325781
$ws.Range("Q6").Value = 0.002017642862888889
$ws.Range("R6").Value = 0.018158785766
$ws.Range("S6").Value = 0.001958801080101452
$ws.Range("T6").Value = 0.001958801080101452

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Icam2"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3556153333333333
$ws.Range("H7").Value = 1.066846
$ws.Range("I7").Value = 0.01586293987910606
$ws.Range("J7").Value = 0.01586293987910605
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04027333333333333
$ws.Range("N7").Value = 0.12082
$ws.Range("O7").Value = 0.8765171465674219
$ws.Range("P7").Value = 0.876517146567422
$ws.Range("Q7").Value = 0.01432181485777778
$ws.Range("R7").Value = 0.12889633372
$ws.Range("S7").Value = 0.01390413879900461
$ws.Range("T7").Value = 0.01390413879900461

# The old sheet additionally had "-> ECs" rows (8:10, MuSCs source)
# that no longer exist in the refreshed output - drop them so the
# sheet dimension shrinks from A1:T10 to A1:T7.
$ws.Rows("8:10").Delete()
